$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.05115733333333
$ws.Range("H2").Value = 126.153472
$ws.Range("I2").Value = 0.1594435451835853
$ws.Range("J2").Value = 0.1594435451835853
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.900405
$ws.Range("N2").Value = 5.701215
$ws.Range("O2").Value = 0.343916381221607
$ws.Range("P2").Value = 0.343916381221607
$ws.Range("Q2").Value = 79.91422965205334
$ws.Range("R2").Value = 719.22806686848
$ws.Range("S2").Value = 0.05483524706868244
$ws.Range("T2").Value = 0.05483524706868245
$ws.Range("G3").Value = 42.05115733333333
$ws.Range("H3").Value = 126.153472
$ws.Range("I3").Value = 0.1594435451835853
$ws.Range("J3").Value = 0.1594435451835853
$ws.Range("O3").Value = 0.3682557063324254
$ws.Range("P3").Value = 0.3682557063324254
$ws.Range("Q3").Value = 85.56984398939021
$ws.Range("R3").Value = 770.1285959045119
$ws.Range("S3").Value = 0.05871599535172719
$ws.Range("T3").Value = 0.0587159953517272
$ws.Range("G4").Value = 42.05115733333333
$ws.Range("H4").Value = 126.153472
$ws.Range("I4").Value = 0.1594435451835853
$ws.Range("J4").Value = 0.1594435451835853
$ws.Range("M4").Value = 1.271690333333333
$ws.Range("N4").Value = 3.815071
$ws.Range("O4").Value = 0.2301378587587904
$ws.Range("P4").Value = 0.2301378587587904
$ws.Range("Q4").Value = 53.4760502862791
$ws.Range("R4").Value = 481.284452576512
$ws.Range("S4").Value = 0.03669399608146077
$ws.Range("T4").Value = 0.03669399608146078
$ws.Range("G5").Value = 42.05115733333333
$ws.Range("H5").Value = 126.153472
$ws.Range("I5").Value = 0.1594435451835853
$ws.Range("J5").Value = 0.1594435451835853
$ws.Range("M5").Value = 0.3187823333333333
$ws.Range("N5").Value = 0.9563470000000001
$ws.Range("O5").Value = 0.05769005368717723
$ws.Range("P5").Value = 0.05769005368717724
$ws.Range("Q5").Value = 13.40516605408711
$ws.Range("R5").Value = 120.646494486784
$ws.Range("S5").Value = 0.009198306681714903
$ws.Range("T5").Value = 0.009198306681714907
$ws.Range("G6").Value = 57.66057933333332
$ws.Range("I6").Value = 0.2186291119973147
$ws.Range("J6").Value = 0.2186291119973148
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.900405
$ws.Range("N6").Value = 5.701215
$ws.Range("O6").Value = 0.343916381221607
$ws.Range("P6").Value = 0.343916381221607
$ws.Range("Q6").Value = 109.5784532679633
$ws.Range("R6").Value = 986.20607941167
$ws.Range("S6").Value = 0.07519013302780991
$ws.Range("T6").Value = 0.07519013302780993
$ws.Range("G7").Value = 57.66057933333332
$ws.Range("I7").Value = 0.2186291119973147
$ws.Range("J7").Value = 0.2186291119973148
$ws.Range("O7").Value = 0.3682557063324254
$ws.Range("P7").Value = 0.3682557063324254
$ws.Range("S7").Value = 0.08051141806340209
$ws.Range("T7").Value = 0.0805114180634021
$ws.Range("G8").Value = 57.66057933333332
$ws.Range("I8").Value = 0.2186291119973147
$ws.Range("J8").Value = 0.2186291119973148
$ws.Range("M8").Value = 1.271690333333333
$ws.Range("N8").Value = 3.815071
$ws.Range("O8").Value = 0.2301378587587904
$ws.Range("P8").Value = 0.2301378587587904
$ws.Range("Q8").Value = 73.32640135259977
$ws.Range("R8").Value = 659.937612173398
$ws.Range("S8").Value = 0.05031483569739779
$ws.Range("T8").Value = 0.0503148356973978
$ws.Range("G9").Value = 57.66057933333332
$ws.Range("I9").Value = 0.2186291119973147
$ws.Range("J9").Value = 0.2186291119973148
$ws.Range("M9").Value = 0.3187823333333333
$ws.Range("N9").Value = 0.9563470000000001
$ws.Range("O9").Value = 0.05769005368717723
$ws.Range("P9").Value = 0.05769005368717724
$ws.Range("Q9").Value = 18.38117402123178
$ws.Range("R9").Value = 165.430566191086
$ws.Range("S9").Value = 0.01261272520870497
$ws.Range("T9").Value = 0.01261272520870497
$ws.Range("G10").Value = 99.15200299999999
$ws.Range("H10").Value = 297.456009
$ws.Range("I10").Value = 0.3759503393701321
$ws.Range("J10").Value = 0.3759503393701321
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.900405
$ws.Range("N10").Value = 5.701215
$ws.Range("O10").Value = 0.343916381221607
$ws.Range("P10").Value = 0.343916381221607
$ws.Range("Q10").Value = 188.428962261215
$ws.Range("R10").Value = 1695.860660350935
$ws.Range("S10").Value = 0.1292954802352109
$ws.Range("T10").Value = 0.1292954802352109
$ws.Range("G11").Value = 99.15200299999999
$ws.Range("H11").Value = 297.456009
$ws.Range("I11").Value = 0.3759503393701321
$ws.Range("J11").Value = 0.3759503393701321
$ws.Range("O11").Value = 0.3682557063324254
$ws.Range("P11").Value = 0.3682557063324254
$ws.Range("Q11").Value = 201.7642787020293
$ws.Range("R11").Value = 1815.878508318264
$ws.Range("S11").Value = 0.138445857770663
$ws.Range("T11").Value = 0.138445857770663
$ws.Range("G12").Value = 99.15200299999999
$ws.Range("H12").Value = 297.456009
$ws.Range("I12").Value = 0.3759503393701321
$ws.Range("J12").Value = 0.3759503393701321
$ws.Range("M12").Value = 1.271690333333333
$ws.Range("N12").Value = 3.815071
$ws.Range("O12").Value = 0.2301378587587904
$ws.Range("P12").Value = 0.2301378587587904
$ws.Range("Q12").Value = 126.0906437457377
$ws.Range("R12").Value = 1134.815793711639
$ws.Range("S12").Value = 0.08652040610228277
$ws.Range("T12").Value = 0.08652040610228279
$ws.Range("G13").Value = 99.15200299999999
$ws.Range("H13").Value = 297.456009
$ws.Range("I13").Value = 0.3759503393701321
$ws.Range("J13").Value = 0.3759503393701321
$ws.Range("M13").Value = 0.3187823333333333
$ws.Range("N13").Value = 0.9563470000000001
$ws.Range("O13").Value = 0.05769005368717723
$ws.Range("P13").Value = 0.05769005368717724
$ws.Range("Q13").Value = 31.60790687101366
$ws.Range("R13").Value = 284.471161839123
$ws.Range("S13").Value = 0.02168859526197542
$ws.Range("T13").Value = 0.02168859526197542
$ws.Range("G14").Value = 64.87322933333333
$ws.Range("H14").Value = 194.619688
$ws.Range("I14").Value = 0.2459770034489679
$ws.Range("J14").Value = 0.2459770034489679
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.900405
$ws.Range("N14").Value = 5.701215
$ws.Range("O14").Value = 0.343916381221607
$ws.Range("P14").Value = 0.343916381221607
$ws.Range("Q14").Value = 123.2854093912133
$ws.Range("R14").Value = 1109.56868452092
$ws.Range("S14").Value = 0.08459552088990377
$ws.Range("T14").Value = 0.08459552088990378
$ws.Range("G15").Value = 64.87322933333333
$ws.Range("H15").Value = 194.619688
$ws.Range("I15").Value = 0.2459770034489679
$ws.Range("J15").Value = 0.2459770034489679
$ws.Range("O15").Value = 0.3682557063324254
$ws.Range("P15").Value = 0.3682557063324254
$ws.Range("Q15").Value = 132.0104478727609
$ws.Range("R15").Value = 1188.094030854848
$ws.Range("S15").Value = 0.09058243514663311
$ws.Range("T15").Value = 0.09058243514663311
$ws.Range("G16").Value = 64.87322933333333
$ws.Range("H16").Value = 194.619688
$ws.Range("I16").Value = 0.2459770034489679
$ws.Range("J16").Value = 0.2459770034489679
$ws.Range("M16").Value = 1.271690333333333
$ws.Range("N16").Value = 3.815071
$ws.Range("O16").Value = 0.2301378587587904
$ws.Range("P16").Value = 0.2301378587587904
$ws.Range("Q16").Value = 82.49865863531645
$ws.Range("R16").Value = 742.487927717848
$ws.Range("S16").Value = 0.05660862087764907
$ws.Range("T16").Value = 0.05660862087764908
$ws.Range("G17").Value = 64.87322933333333
$ws.Range("H17").Value = 194.619688
$ws.Range("I17").Value = 0.2459770034489679
$ws.Range("J17").Value = 0.2459770034489679
$ws.Range("M17").Value = 0.3187823333333333
$ws.Range("N17").Value = 0.9563470000000001
$ws.Range("O17").Value = 0.05769005368717723
$ws.Range("P17").Value = 0.05769005368717724
$ws.Range("Q17").Value = 20.68043941774844
$ws.Range("R17").Value = 186.123954759736
$ws.Range("S17").Value = 0.01419042653478193
$ws.Range("T17").Value = 0.01419042653478194
